# Update ELMo-BiLSTM training history rows (3, 5, 8) with new result values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$rowRange3 = $ws.Range("E3:BB3")
$rowRange3.NumberFormat = "@"
$ws.Range("E3").Value = "0.761827"
$ws.Range("F3").Value = "0.832407"
$ws.Range("G3").Value = "0.845679"
$ws.Range("H3").Value = "0.831776"
$ws.Range("I3").Value = "0.820312"
$ws.Range("J3").Value = "0.852334"
$ws.Range("K3").Value = "0.852459"
$ws.Range("L3").Value = "0.846448"
$ws.Range("M3").Value = "0.854054"
$ws.Range("N3").Value = "0.819856"
$ws.Range("O3").Value = "0.859756"
$ws.Range("P3").Value = "0.835218"
$ws.Range("Q3").Value = "0.877301"
$ws.Range("R3").Value = "0.869166"
$ws.Range("S3").Value = "0.811321"
$ws.Range("T3").Value = "0.879509"
$ws.Range("U3").Value = "0.824645"
$ws.Range("V3").Value = "0.840804"
$ws.Range("W3").Value = "0.828927"
$ws.Range("X3").Value = "0.8382"
$ws.Range("Y3").Value = "0.88408"
$ws.Range("Z3").Value = "0.851562"
$ws.Range("AA3").Value = "0.86262"
$ws.Range("AB3").Value = "0.852433"
$ws.Range("AC3").Value = "0.848012"
$ws.Range("AD3").Value = "0.848954"
$ws.Range("AE3").Value = "0.89112"
$ws.Range("AF3").Value = "0.864738"
$ws.Range("AG3").Value = "0.830964"
$ws.Range("AH3").Value = "0.860876"
$ws.Range("AI3").Value = "0.864101"
$ws.Range("AJ3").Value = "0.84492"
$ws.Range("AK3").Value = "0.818959"
$ws.Range("AL3").Value = "0.87404"
$ws.Range("AM3").Value = "0.837713"
$ws.Range("AN3").Value = "0.829157"
$ws.Range("AO3").Value = "0.86387"
$ws.Range("AP3").Value = "0.838168"
$ws.Range("AQ3").Value = "0.868827"
$ws.Range("AR3").Value = "0.879087"
$ws.Range("AS3").Value = "0.850467"
$ws.Range("AT3").Value = "0.892664"
$ws.Range("AU3").Value = "0.854003"
$ws.Range("AV3").Value = "0.839062"
$ws.Range("AW3").Value = "0.850965"
$ws.Range("AX3").Value = "0.837817"
$ws.Range("AY3").Value = "0.868217"
$ws.Range("AZ3").Value = "0.837817"
$ws.Range("BA3").Value = "0.848485"
$ws.Range("BB3").Value = "0.855807"
$rowRange3.Style = "Normal"

# Row 5
$rowRange5 = $ws.Range("E5:BB5")
$rowRange5.NumberFormat = "@"
$ws.Range("E5").Value = "0.707885"
$ws.Range("F5").Value = "0.81157"
$ws.Range("G5").Value = "0.8382"
$ws.Range("H5").Value = "0.800621"
$ws.Range("I5").Value = "0.837061"
$ws.Range("J5").Value = "0.84441"
$ws.Range("K5").Value = "0.849273"
$ws.Range("L5").Value = "0.83007"
$ws.Range("M5").Value = "0.827103"
$ws.Range("N5").Value = "0.816746"
$ws.Range("O5").Value = "0.827423"
$ws.Range("P5").Value = "0.807018"
$ws.Range("Q5").Value = "0.827533"
$ws.Range("R5").Value = "0.865533"
$ws.Range("S5").Value = "0.865194"
$ws.Range("T5").Value = "0.832162"
$ws.Range("U5").Value = "0.840376"
$ws.Range("V5").Value = "0.86537"
$ws.Range("W5").Value = "0.841521"
$ws.Range("X5").Value = "0.849336"
$ws.Range("Y5").Value = "0.855129"
$ws.Range("Z5").Value = "0.846395"
$ws.Range("AA5").Value = "0.840874"
$ws.Range("AB5").Value = "0.866302"
$ws.Range("AC5").Value = "0.870769"
$ws.Range("AD5").Value = "0.86"
$ws.Range("AE5").Value = "0.850898"
$ws.Range("AF5").Value = "0.861611"
$ws.Range("AG5").Value = "0.887865"
$ws.Range("AH5").Value = "0.860558"
$ws.Range("AI5").Value = "0.85446"
$ws.Range("AJ5").Value = "0.849765"
$ws.Range("AK5").Value = "0.886662"
$ws.Range("AL5").Value = "0.862042"
$ws.Range("AM5").Value = "0.863287"
$ws.Range("AN5").Value = "0.853168"
$ws.Range("AO5").Value = "0.856474"
$ws.Range("AP5").Value = "0.828235"
$ws.Range("AQ5").Value = "0.851424"
$ws.Range("AR5").Value = "0.843505"
$ws.Range("AS5").Value = "0.819672"
$ws.Range("AT5").Value = "0.863636"
$ws.Range("AU5").Value = "0.844136"
$ws.Range("AV5").Value = "0.860215"
$ws.Range("AW5").Value = "0.85559"
$ws.Range("AX5").Value = "0.853583"
$ws.Range("AY5").Value = "0.842022"
$ws.Range("AZ5").Value = "0.870039"
$ws.Range("BA5").Value = "0.857363"
$ws.Range("BB5").Value = "0.843336"
$rowRange5.Style = "Normal"

# Row 8
$rowRange8 = $ws.Range("E8:BB8")
$rowRange8.NumberFormat = "@"
$ws.Range("E8").Value = "0.806478"
$ws.Range("F8").Value = "0.821681"
$ws.Range("G8").Value = "0.870241"
$ws.Range("H8").Value = "0.869362"
$ws.Range("I8").Value = "0.800959"
$ws.Range("J8").Value = "0.864615"
$ws.Range("K8").Value = "0.848243"
$ws.Range("L8").Value = "0.835752"
$ws.Range("M8").Value = "0.82825"
$ws.Range("N8").Value = "0.847405"
$ws.Range("O8").Value = "0.874214"
$ws.Range("P8").Value = "0.889746"
$ws.Range("Q8").Value = "0.834217"
$ws.Range("R8").Value = "0.897554"
$ws.Range("S8").Value = "0.856705"
$ws.Range("T8").Value = "0.809339"
$ws.Range("U8").Value = "0.863741"
$ws.Range("V8").Value = "0.851163"
$ws.Range("W8").Value = "0.83728"
$ws.Range("X8").Value = "0.844479"
$ws.Range("Y8").Value = "0.869366"
$ws.Range("Z8").Value = "0.83836"
$ws.Range("AA8").Value = "0.840273"
$ws.Range("AB8").Value = "0.881098"
$ws.Range("AC8").Value = "0.83675"
$ws.Range("AD8").Value = "0.831128"
$ws.Range("AE8").Value = "0.860856"
$ws.Range("AF8").Value = "0.852053"
$ws.Range("AG8").Value = "0.848302"
$ws.Range("AH8").Value = "0.864615"
$ws.Range("AI8").Value = "0.862805"
$ws.Range("AJ8").Value = "0.880551"
$ws.Range("AK8").Value = "0.885978"
$ws.Range("AL8").Value = "0.833856"
$ws.Range("AM8").Value = "0.854216"
$ws.Range("AN8").Value = "0.872837"
$ws.Range("AO8").Value = "0.873239"
$ws.Range("AP8").Value = "0.86296"
$ws.Range("AQ8").Value = "0.86758"
$ws.Range("AR8").Value = "0.853792"
$ws.Range("AS8").Value = "0.862654"
$ws.Range("AT8").Value = "0.84488"
$ws.Range("AU8").Value = "0.860429"
$ws.Range("AV8").Value = "0.882353"
$ws.Range("AW8").Value = "0.850115"
$ws.Range("AX8").Value = "0.865109"
$ws.Range("AY8").Value = "0.893553"
$ws.Range("AZ8").Value = "0.881024"
$ws.Range("BA8").Value = "0.844815"
$ws.Range("BB8").Value = "0.874433"
$rowRange8.Style = "Normal"
